$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $newVal) {
    $cell = $ws.Range($rangeAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $newVal
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '63.796.87'
$ws.Range('E2').Value = '  -0.17%  '
Set-TextValue 'D3' '2.733.53'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '564.56'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('E7').Value = '  +0.04%  '
Set-TextValue 'D8' '0.595'
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('E10').Value = '  +4.34%  '
Set-TextValue 'D11' '5.61'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').Value = '  -1.89%  '
Set-TextValue 'D13' '3.220.02'
$ws.Range('E13').Value = '  -0.57%  '
Set-TextValue 'D14' '26.89'
$ws.Range('E14').Value = '  +1.45%  '
Set-TextValue 'D15' '63.642.09'
$ws.Range('E15').Value = '  +0.17%  '
Set-TextValue 'D16' '0.0000150'
$ws.Range('E16').Value = '  -0.56%  '
Set-TextValue 'D17' '2.738.70'
$ws.Range('E17').Value = '  -0.64%  '
Set-TextValue 'D18' '12.28'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('E19').Value = '  -1.81%  '
Set-TextValue 'D20' '355.44'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('E21').Value = '  -1.63%  '
Set-TextValue 'D22' '0.998'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E23').Value = '  -3.00%  '
Set-TextValue 'D24' '64.28'
$ws.Range('E24').Value = '  -1.32%  '
Set-TextValue 'D25' '0.170'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('E26').Value = '  -0.07%  '
Set-TextValue 'D27' '8.35'
$ws.Range('E27').Value = '  -1.36%  '
Set-TextValue 'D28' '0.0₃0908'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('E29').Value = '  +2.70%  '
$ws.Range('E30').Value = '  +8.48%  '
Set-TextValue 'D31' '7.14'
$ws.Range('E31').Value = '  +1.18%  '
Set-TextValue 'D32' '167.01'
$ws.Range('E32').Value = '  -1.05%  '
$ws.Range('E33').Value = '  +0.14%  '
Set-TextValue 'D34' '20.03'
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('E35').Value = '  +1.98%  '
Set-TextValue 'D36' '0.998'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  +0.58%  '
Set-TextValue 'D38' '0.972'
$ws.Range('E38').Value = '  -1.47%  '
Set-TextValue 'D39' '346.24'
$ws.Range('E39').Value = '  +4.38%  '
Set-TextValue 'D40' '6.28'
$ws.Range('E40').Value = '  +1.88%  '
Set-TextValue 'D41' '4.07'
$ws.Range('E41').Value = '  -1.75%  '
Set-TextValue 'D42' '38.64'
$ws.Range('E42').Value = '  -0.88%  '
Set-TextValue 'D43' '21.75'
$ws.Range('E43').Value = '  +1.17%  '
Set-TextValue 'D44' '21.09'
$ws.Range('E44').Value = '  -1.70%  '
Set-TextValue 'D45' '0.0581'
$ws.Range('E45').Value = '  -0.93%  '
Set-TextValue 'D46' '0.631'
$ws.Range('E46').Value = '  +0.96%  '
$ws.Range('E47').Value = '  -1.74%  '
Set-TextValue 'D48' '0.0994'
$ws.Range('E48').Value = '  -1.30%  '
Set-TextValue 'D49' '132.38'
$ws.Range('E49').Value = '  -1.98%  '
Set-TextValue 'D51' '11.05'
$ws.Range('E51').Value = '  +0.22%  '
